# Apply the "more items, monsters and affixes" update to monsters.xlsx
#
# Summary of the change:
#   - "Monsters" sheet: 8 new monster rows appended (rows 23-30), and
#     column N (attack_range) widened from 15 to 16 chars.
#   - "Monsters Skills" sheet: 16 new rows appended (rows 44-59), two
#     (Accuracy/Dodge) for each of the 8 new monsters.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Monsters"
$ws2 = $wb.Worksheets.Item(2)   # "Monsters Skills"

# --- Monsters sheet (sheet1): add rows 23-30 ---
# Row 23: Labyrinth Fiend
$ws1.Cells.Item(23, 1).Value = "Labyrinth Fiend"
$ws1.Cells.Item(23, 2).Value = 750
$ws1.Cells.Item(23, 3).Value = 730
$ws1.Cells.Item(23, 4).Value = 740
$ws1.Cells.Item(23, 5).Value = 755
$ws1.Cells.Item(23, 6).Value = 750
$ws1.Cells.Item(23, 7).Value = 140
$ws1.Cells.Item(23, 8).Value = 750
$ws1.Cells.Item(23, 9).Value = "chr"
$ws1.Cells.Item(23, 10).Value = 10
$ws1.Cells.Item(23, 11).Value = 0.01
$ws1.Cells.Item(23, 12).Value = 2000
$ws1.Cells.Item(23, 13).Value = "12670-13800"
$ws1.Cells.Item(23, 14).Value = "42000-49000"
$ws1.Cells.Item(23, 15).Value = 1
$ws1.Cells.Item(23, 16).Value = "Key To The Labyrinth"
$ws1.Cells.Item(23, 17).Value = 0.01
$ws1.Cells.Item(23, 18).Value = "Surface"

# Row 24: Deaths Horsemen
$ws1.Cells.Item(24, 1).Value = "Deaths Horsemen"
$ws1.Cells.Item(24, 2).Value = 720
$ws1.Cells.Item(24, 3).Value = 780
$ws1.Cells.Item(24, 4).Value = 760
$ws1.Cells.Item(24, 5).Value = 750
$ws1.Cells.Item(24, 6).Value = 785
$ws1.Cells.Item(24, 7).Value = 145
$ws1.Cells.Item(24, 8).Value = 790
$ws1.Cells.Item(24, 9).Value = "int"
$ws1.Cells.Item(24, 10).Value = 10
$ws1.Cells.Item(24, 11).Value = 0.04
$ws1.Cells.Item(24, 12).Value = 2300
$ws1.Cells.Item(24, 13).Value = "14000-15000"
$ws1.Cells.Item(24, 14).Value = "53000-58000"
$ws1.Cells.Item(24, 15).Value = 1
$ws1.Cells.Item(24, 18).Value = "Surface"

# Row 25: Celestial Entity
$ws1.Cells.Item(25, 1).Value = "Celestial Entity"
$ws1.Cells.Item(25, 2).Value = 830
$ws1.Cells.Item(25, 3).Value = 845
$ws1.Cells.Item(25, 4).Value = 830
$ws1.Cells.Item(25, 5).Value = 840
$ws1.Cells.Item(25, 6).Value = 840
$ws1.Cells.Item(25, 7).Value = 155
$ws1.Cells.Item(25, 8).Value = 840
$ws1.Cells.Item(25, 9).Value = "dur"
$ws1.Cells.Item(25, 10).Value = 12
$ws1.Cells.Item(25, 11).Value = 0.02
$ws1.Cells.Item(25, 12).Value = 2450
$ws1.Cells.Item(25, 13).Value = "15000-16500"
$ws1.Cells.Item(25, 14).Value = "62370-73000"
$ws1.Cells.Item(25, 15).Value = 1
$ws1.Cells.Item(25, 18).Value = "Surface"

# Row 26: Shade of Light
$ws1.Cells.Item(26, 1).Value = "Shade of Light"
$ws1.Cells.Item(26, 2).Value = 850
$ws1.Cells.Item(26, 3).Value = 850
$ws1.Cells.Item(26, 4).Value = 860
$ws1.Cells.Item(26, 5).Value = 840
$ws1.Cells.Item(26, 6).Value = 830
$ws1.Cells.Item(26, 7).Value = 160
$ws1.Cells.Item(26, 8).Value = 870
$ws1.Cells.Item(26, 9).Value = "dex"
$ws1.Cells.Item(26, 10).Value = 12
$ws1.Cells.Item(26, 11).Value = 0.05
$ws1.Cells.Item(26, 12).Value = 2600
$ws1.Cells.Item(26, 13).Value = "17000-20000"
$ws1.Cells.Item(26, 14).Value = "75000-80000"
$ws1.Cells.Item(26, 15).Value = 1
$ws1.Cells.Item(26, 18).Value = "Surface"

# Row 27: Ruby Fiend
$ws1.Cells.Item(27, 1).Value = "Ruby Fiend"
$ws1.Cells.Item(27, 2).Value = 920
$ws1.Cells.Item(27, 3).Value = 890
$ws1.Cells.Item(27, 4).Value = 900
$ws1.Cells.Item(27, 5).Value = 880
$ws1.Cells.Item(27, 6).Value = 870
$ws1.Cells.Item(27, 7).Value = 200
$ws1.Cells.Item(27, 8).Value = 920
$ws1.Cells.Item(27, 9).Value = "str"
$ws1.Cells.Item(27, 10).Value = 12
$ws1.Cells.Item(27, 11).Value = 0.01
$ws1.Cells.Item(27, 12).Value = 3000
$ws1.Cells.Item(27, 13).Value = "22000-24000"
$ws1.Cells.Item(27, 14).Value = "82000-86000"
$ws1.Cells.Item(27, 15).Value = 1
$ws1.Cells.Item(27, 18).Value = "Surface"

# Row 28: Astral Hell God
$ws1.Cells.Item(28, 1).Value = "Astral Hell God"
$ws1.Cells.Item(28, 2).Value = 950
$ws1.Cells.Item(28, 3).Value = 950
$ws1.Cells.Item(28, 4).Value = 940
$ws1.Cells.Item(28, 5).Value = 950
$ws1.Cells.Item(28, 6).Value = 930
$ws1.Cells.Item(28, 7).Value = 230
$ws1.Cells.Item(28, 8).Value = 950
$ws1.Cells.Item(28, 9).Value = "chr"
$ws1.Cells.Item(28, 10).Value = 15
$ws1.Cells.Item(28, 11).Value = 0.02
$ws1.Cells.Item(28, 12).Value = 3300
$ws1.Cells.Item(28, 13).Value = "25000-27000"
$ws1.Cells.Item(28, 14).Value = "90000-100000"
$ws1.Cells.Item(28, 15).Value = 1
$ws1.Cells.Item(28, 18).Value = "Surface"

# Row 29: Jester Of Beleth
$ws1.Cells.Item(29, 1).Value = "Jester Of Beleth"
$ws1.Cells.Item(29, 2).Value = 960
$ws1.Cells.Item(29, 3).Value = 970
$ws1.Cells.Item(29, 4).Value = 950
$ws1.Cells.Item(29, 5).Value = 975
$ws1.Cells.Item(29, 6).Value = 978
$ws1.Cells.Item(29, 7).Value = 260
$ws1.Cells.Item(29, 8).Value = 970
$ws1.Cells.Item(29, 9).Value = "int"
$ws1.Cells.Item(29, 10).Value = 15
$ws1.Cells.Item(29, 11).Value = 0.04
$ws1.Cells.Item(29, 12).Value = 3400
$ws1.Cells.Item(29, 13).Value = "28000-29500"
$ws1.Cells.Item(29, 14).Value = "110000-130000"
$ws1.Cells.Item(29, 15).Value = 1
$ws1.Cells.Item(29, 18).Value = "Surface"

# Row 30: Satanic Cult Leader
$ws1.Cells.Item(30, 1).Value = "Satanic Cult Leader"
$ws1.Cells.Item(30, 2).Value = 999
$ws1.Cells.Item(30, 3).Value = 999
$ws1.Cells.Item(30, 4).Value = 999
$ws1.Cells.Item(30, 5).Value = 999
$ws1.Cells.Item(30, 6).Value = 999
$ws1.Cells.Item(30, 7).Value = 300
$ws1.Cells.Item(30, 8).Value = 100000
$ws1.Cells.Item(30, 9).Value = "str"
$ws1.Cells.Item(30, 10).Value = 20
$ws1.Cells.Item(30, 11).Value = 0.01
$ws1.Cells.Item(30, 12).Value = 5000
$ws1.Cells.Item(30, 13).Value = "30000-50000"
$ws1.Cells.Item(30, 14).Value = "150000-200000"
$ws1.Cells.Item(30, 18).Value = "Surface"

# Widen column N (attack_range) from 15 to 16 characters of bestFit width.
# ColumnWidth uses a slightly different unit than the raw OOXML "width"
# attribute (offset by the default cell padding), so compensate to land
# exactly on width=16 in the saved file.
$ws1.Columns.Item(14).ColumnWidth = 15.166666666666666

# --- Monsters Skills sheet (sheet2): add rows 44-59 ---
# Two rows (Accuracy, Dodge) per new monster, in the same order they were
# added to the Monsters sheet.
$ws2.Cells.Item(44, 1).Value = "Labyrinth Fiend"
$ws2.Cells.Item(44, 3).Value = 39
$ws2.Cells.Item(44, 4).Value = 0
$ws2.Cells.Item(44, 7).Value = "Accuracy"

$ws2.Cells.Item(45, 1).Value = "Labyrinth Fiend"
$ws2.Cells.Item(45, 3).Value = 0
$ws2.Cells.Item(45, 4).Value = 0
$ws2.Cells.Item(45, 7).Value = "Dodge"

$ws2.Cells.Item(46, 1).Value = "Deaths Horsemen"
$ws2.Cells.Item(46, 3).Value = 40
$ws2.Cells.Item(46, 4).Value = 0
$ws2.Cells.Item(46, 7).Value = "Accuracy"

$ws2.Cells.Item(47, 1).Value = "Deaths Horsemen"
$ws2.Cells.Item(47, 3).Value = 35
$ws2.Cells.Item(47, 4).Value = 0
$ws2.Cells.Item(47, 7).Value = "Dodge"

$ws2.Cells.Item(48, 1).Value = "Celestial Entity"
$ws2.Cells.Item(48, 3).Value = 42
$ws2.Cells.Item(48, 4).Value = 0
$ws2.Cells.Item(48, 7).Value = "Accuracy"

$ws2.Cells.Item(49, 1).Value = "Celestial Entity"
$ws2.Cells.Item(49, 3).Value = 38
$ws2.Cells.Item(49, 4).Value = 0
$ws2.Cells.Item(49, 7).Value = "Dodge"

$ws2.Cells.Item(50, 1).Value = "Shade of Light"
$ws2.Cells.Item(50, 3).Value = 43
$ws2.Cells.Item(50, 4).Value = 0
$ws2.Cells.Item(50, 7).Value = "Accuracy"

$ws2.Cells.Item(51, 1).Value = "Shade of Light"
$ws2.Cells.Item(51, 3).Value = 45
$ws2.Cells.Item(51, 4).Value = 0
$ws2.Cells.Item(51, 7).Value = "Dodge"

$ws2.Cells.Item(52, 1).Value = "Ruby Fiend"
$ws2.Cells.Item(52, 3).Value = 45
$ws2.Cells.Item(52, 4).Value = 0
$ws2.Cells.Item(52, 7).Value = "Accuracy"

$ws2.Cells.Item(53, 1).Value = "Ruby Fiend"
$ws2.Cells.Item(53, 3).Value = 42
$ws2.Cells.Item(53, 4).Value = 0
$ws2.Cells.Item(53, 7).Value = "Dodge"

$ws2.Cells.Item(54, 1).Value = "Astral Hell God"
$ws2.Cells.Item(54, 3).Value = 45
$ws2.Cells.Item(54, 4).Value = 0
$ws2.Cells.Item(54, 7).Value = "Accuracy"

$ws2.Cells.Item(55, 1).Value = "Astral Hell God"
$ws2.Cells.Item(55, 3).Value = 45
$ws2.Cells.Item(55, 4).Value = 0
$ws2.Cells.Item(55, 7).Value = "Dodge"

$ws2.Cells.Item(56, 1).Value = "Jester Of Beleth"
$ws2.Cells.Item(56, 3).Value = 45
$ws2.Cells.Item(56, 4).Value = 0
$ws2.Cells.Item(56, 7).Value = "Accuracy"

$ws2.Cells.Item(57, 1).Value = "Jester Of Beleth"
$ws2.Cells.Item(57, 3).Value = 45
$ws2.Cells.Item(57, 4).Value = 0
$ws2.Cells.Item(57, 7).Value = "Dodge"

$ws2.Cells.Item(58, 1).Value = "Satanic Cult Leader"
$ws2.Cells.Item(58, 3).Value = 45
$ws2.Cells.Item(58, 4).Value = 0
$ws2.Cells.Item(58, 7).Value = "Accuracy"

$ws2.Cells.Item(59, 1).Value = "Satanic Cult Leader"
$ws2.Cells.Item(59, 3).Value = 45
$ws2.Cells.Item(59, 4).Value = 0
$ws2.Cells.Item(59, 7).Value = "Dodge"
